$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Seed row 8 with the formatting from the row above (A7:L7) so the new
# date cell picks up the existing date style instead of a brand-new one.
$ws.Range("A7:L7").Copy($ws.Range("A8:L8"))

# Append the new weekly data row (week of Dec 9, 2020 -> serial 44174),
# extending the CasesByAge table from row 7 to row 8.
$ws.Range("A8").Value = 44174
$ws.Range("B8").Value = 8201
$ws.Range("C8").Value = 9746
$ws.Range("D8").Value = 8174
$ws.Range("E8").Value = 6835
$ws.Range("F8").Value = 7141
$ws.Range("G8").Value = 4683
$ws.Range("H8").Value = 2317
$ws.Range("I8").Value = 1676
$ws.Range("J8").Value = 69
$ws.Range("K8").Value = 82
$ws.Range("L8").Value = 50
